$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case row 18 mirrors the formatting of row 16 (same style pattern:
# s=2, s=4, s=6, s=2, s=2 with a tall wrapped-text description cell), so
# copy formats from row 16 down into row 18 first.
$ws.Range("A16:E16").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)

# Fill in the new RCC114 test case. Writing A (TCID) then C (Description)
# then B (Jira id) then D (Runmode) so new shared-string entries land in
# the same order as the target workbook.
$ws.Range("A18").Value = "RCC114"
$ws.Range("C18").Value = "Verify that user is able to add an article to the multiple groups from search results page.||Verify that user is able to add a post to the  multiple groups from Search results page.||Verify that user is able to add a patent to the  multiple groups from search results page."
$ws.Range("B18").Value = "OPQA-3468||OPQA-3472||OPQA-3476"
$ws.Range("D18").Value = "Y"

# Row height for the new row (wrapped 3-line description).
$ws.Rows.Item(18).RowHeight = 45

# Move the active selection down past the newly-added row.
$null = $ws.Range("C19").Select()
